$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.459.60"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "1.828.73"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").Value = "'330.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").Value = "'0.4587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "'0.3828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").Value = "'46.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.07914"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.9693"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "'21.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "1.867.65"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'5.872"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'7.048"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'88.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'0.06652"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'0.00001030"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "'17.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "27.458.90"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "'5.335"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "'10.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "2.041.42"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "'157.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "'19.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").Value = "'2.059"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "'5.233"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'118.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").Value = "'0.9481"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'0.09284"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "'3.563"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "'5.237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "'1.315"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'0.02203"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.05923"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.154"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").Value = "'0.5779"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value = "'0.1839"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'10.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "'1.288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "'0.5480"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").Value = "'12.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'1.863"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "'0.06648"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'109.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "'1.039"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "'1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
